$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 16; $r++) {
    # birth_year column (E) -> 1963
    $ws.Cells.Item($r, 5).Value2 = 1963
    # age column (G) -> increment existing value by 1
    $currentAge = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 7).Value2 = $currentAge + 1
}
